# Auto-generated edit script: update scraped market-price derived values
# in the Zodiark_Profits workbook (one worksheet per crafting job table).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4711
$ws.Range("J43").Value = 2996
$ws.Range("L43").Value = 2996
$ws.Range("N43").Value = -3134
$ws.Range("H80").Value = 477.46155
$ws.Range("J80").Value = 737.1539
$ws.Range("L80").Value = 2211.4617
$ws.Range("N80").Value = -4207.4617
$ws.Range("H82").Value = 1527.6666
$ws.Range("I82").Value = 1484.8182
$ws.Range("K82").Value = 4454.4546
$ws.Range("M82").Value = -4048.4546
$ws.Range("H83").Value = 477.46155
$ws.Range("J83").Value = 737.1539
$ws.Range("L83").Value = 6634.3851
$ws.Range("N83").Value = -16618.3851
$ws.Range("H85").Value = 1527.6666
$ws.Range("I85").Value = 1484.8182
$ws.Range("K85").Value = 4454.4546
$ws.Range("M85").Value = -3050.4546
$ws.Range("H132").Value = 1517.4375
$ws.Range("I132").Value = 1517.4375
$ws.Range("K132").Value = 4552.3125
$ws.Range("M132").Value = -2022.3125
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 6900
$ws.Range("I135").Value = 6900
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 62100
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -59565
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 72998
$ws.Range("J136").Value = 72998
$ws.Range("L136").Value = 72998
$ws.Range("N136").Value = -83198
$ws.Range("H137").Value = 25002596
$ws.Range("I137").Value = 31252848
$ws.Range("K137").Value = 93758544
$ws.Range("M137").Value = -93755994
$ws.Range("H138").Value = 1981.4906
$ws.Range("J138").Value = 1974.84
$ws.Range("L138").Value = 5924.52
$ws.Range("N138").Value = -16204.52
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1899.5186
$ws.Range("I88").Value = 1138.2222
$ws.Range("K88").Value = 1138.2222
$ws.Range("M88").Value = -732.2221999999999
$ws.Range("H91").Value = 1899.5186
$ws.Range("I91").Value = 1138.2222
$ws.Range("K91").Value = 1138.2222
$ws.Range("M91").Value = 265.7778000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3040.1667
$ws.Range("I20").Value = 1030.3334
$ws.Range("J20").Value = 5050
$ws.Range("K20").Value = 1030.3334
$ws.Range("L20").Value = 5050
$ws.Range("M20").Value = -783.3334
$ws.Range("N20").Value = -5544
$ws.Range("H86").Value = 3528.4
$ws.Range("I86").Value = 3216.1667
$ws.Range("K86").Value = 3216.1667
$ws.Range("M86").Value = -2093.1667
$ws.Range("H89").Value = 3528.4
$ws.Range("I89").Value = 3216.1667
$ws.Range("K89").Value = 16080.8335
$ws.Range("M89").Value = -10464.8335
$ws.Range("H107").Value = 6116.3335
$ws.Range("I107").Value = 6135
$ws.Range("K107").Value = 6135
$ws.Range("M107").Value = -4215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6538.25
$ws.Range("I16").Value = 8755.25
$ws.Range("J16").Value = 4321.25
$ws.Range("K16").Value = 8755.25
$ws.Range("L16").Value = 4321.25
$ws.Range("M16").Value = -8468.25
$ws.Range("N16").Value = -4895.25
$ws.Range("H35").Value = 278.2353
$ws.Range("I35").Value = 278.2353
$ws.Range("K35").Value = 278.2353
$ws.Range("M35").Value = 15.7647
$ws.Range("H54").Value = 16000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 16000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 16000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -17316
$ws.Range("H58").Value = 4755.2104
$ws.Range("I58").Value = 1989.5714
$ws.Range("K58").Value = 1989.5714
$ws.Range("M58").Value = -1786.5714
$ws.Range("H107").Value = 823
$ws.Range("I107").Value = 454.1111
$ws.Range("K107").Value = 454.1111
$ws.Range("M107").Value = 1465.8889
$ws.Range("H113").Value = 6538.25
$ws.Range("I113").Value = 8755.25
$ws.Range("J113").Value = 4321.25
$ws.Range("K113").Value = 8755.25
$ws.Range("L113").Value = 4321.25
$ws.Range("M113").Value = -6585.25
$ws.Range("N113").Value = -8661.25
$ws.Range("H136").Value = 4755.2104
$ws.Range("I136").Value = 1989.5714
$ws.Range("K136").Value = 5968.7142
$ws.Range("M136").Value = -3418.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 7789.9
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 8544.333000000001
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 25632.999
$ws.Range("M9").Value = -2776
$ws.Range("N9").Value = -26080.999
$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 500
$ws.Range("K68").Value = 1500
$ws.Range("M68").Value = -689
$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 500
$ws.Range("K71").Value = 4500
$ws.Range("M71").Value = -444
$ws.Range("H80").Value = 3249.5
$ws.Range("H83").Value = 3249.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H52").Value = 44733
$ws.Range("J52").Value = 44733
$ws.Range("L52").Value = 44733
$ws.Range("N52").Value = -45251
$ws.Range("H132").Value = 3234.0625
$ws.Range("I132").Value = 3082.9285
$ws.Range("J132").Value = 4292
$ws.Range("K132").Value = 9248.7855
$ws.Range("L132").Value = 12876
$ws.Range("M132").Value = -6718.7855
$ws.Range("N132").Value = -17936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1670.4736
$ws.Range("J82").Value = 3984
$ws.Range("L82").Value = 3984
$ws.Range("N82").Value = -4706
$ws.Range("H85").Value = 1670.4736
$ws.Range("J85").Value = 3984
$ws.Range("L85").Value = 3984
$ws.Range("N85").Value = -6480

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5308.25
$ws.Range("J81").Value = 6344.3335
$ws.Range("L81").Value = 12688.667
$ws.Range("N81").Value = -14810.667
$ws.Range("H84").Value = 5308.25
$ws.Range("J84").Value = 6344.3335
$ws.Range("L84").Value = 63443.335
$ws.Range("N84").Value = -74051.33499999999
$ws.Range("H95").Value = 1000000000
$ws.Range("J95").Value = 1000000000
$ws.Range("L95").Value = 1000000000
$ws.Range("N95").Value = -1000005492
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 1029.8572
$ws.Range("I107").Value = 902.25
$ws.Range("K107").Value = 2706.75
$ws.Range("M107").Value = -786.75
$ws.Range("H136").Value = 3201.5356
$ws.Range("I136").Value = 2190.9092
$ws.Range("J136").Value = 6907.1665
$ws.Range("K136").Value = 6572.7276
$ws.Range("L136").Value = 20721.4995
$ws.Range("M136").Value = -4022.7276
$ws.Range("N136").Value = -25821.4995

